$d = $word.ActiveDocument

# Locate the paragraph that contains the "4.0-Basic-User-Guide" hyperlink --
# the anchor point after which the new content must be inserted.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*4.0-Basic-User-Guide*") {
        $anchorIndex = $i
        break
    }
}

$rng = $d.Paragraphs.Item($anchorIndex).Range
$rng.Collapse(0)

# Insert the first new paragraph -- this becomes the blank "<w:p/>" from the
# diff, so make sure it stays completely empty (no stray run).
$rng.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Item($anchorIndex + 1)
$blankRng = $blankPara.Range
$xmlBlank = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$blankRng.InsertXML($xmlBlank)

# Insert a second new paragraph right after the blank one; this will hold the
# "code . --no-sandbox --disable-gpu-sandbox" content, reproducing the runs
# and proofing-error markers exactly as they appear in the target document.
$blankPara = $d.Paragraphs.Item($anchorIndex + 1)
$blankRng = $blankPara.Range
$blankRng.InsertParagraphAfter()

$codePara = $d.Paragraphs.Item($anchorIndex + 2)
$codeRng = $codePara.Range

$xmlCode = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">code </w:t></w:r><w:r><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>--no-sandbox --disable-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gpu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-sandbox</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$codeRng.InsertXML($xmlCode)
